$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value2 = "58.378.88"
$ws.Cells.Item(2,5).Value2 = "  -3.29%  "

# Row 3
$ws.Cells.Item(3,4).Value2 = "2.984.40"
$ws.Cells.Item(3,5).Value2 = "  +0.20%  "

# Row 4
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value2 = "0.997"
$ws.Cells.Item(4,5).Value2 = "  -0.33%  "

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value2 = "559.64"
$ws.Cells.Item(5,5).Value2 = "  -1.30%  "

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value2 = "130.85"
$ws.Cells.Item(6,5).Value2 = "  +5.44%  "

# Row 7
$ws.Cells.Item(7,5).Value2 = "  -0.08%  "

# Row 8
$ws.Cells.Item(8,2).Value2 = "LidoStakedEther"
$ws.Cells.Item(8,3).Value2 = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Cells.Item(8,4).Value2 = "2.970.89"
$ws.Cells.Item(8,5).Value2 = "  -0.16%  "

# Row 9
$ws.Cells.Item(9,2).Value2 = "XRP"
$ws.Cells.Item(9,3).Value2 = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value2 = "0.513"
$ws.Cells.Item(9,5).Value2 = "  +3.26%  "

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value2 = "0.130"
$ws.Cells.Item(10,5).Value2 = "  -1.60%  "

# Row 11
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value2 = "4.87"
$ws.Cells.Item(11,5).Value2 = "  -4.99%  "

# Row 12
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value2 = "0.450"
$ws.Cells.Item(12,5).Value2 = "  +3.58%  "

# Row 13
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value2 = "0.0000223"
$ws.Cells.Item(13,5).Value2 = "  +0.67%  "

# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value2 = "33.20"
$ws.Cells.Item(14,5).Value2 = "  +2.40%  "

# Row 15
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value2 = "0.121"
$ws.Cells.Item(15,5).Value2 = "  +1.92%  "

# Row 16
$ws.Cells.Item(16,4).Value2 = "3.450.52"
$ws.Cells.Item(16,5).Value2 = "  -0.38%  "

# Row 17
$ws.Cells.Item(17,4).Value2 = "2.960.12"
$ws.Cells.Item(17,5).Value2 = "  -0.44%  "

# Row 18
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value2 = "6.75"
$ws.Cells.Item(18,5).Value2 = "  +10.11%  "

# Row 19
$ws.Cells.Item(19,4).Value2 = "58.005.85"
$ws.Cells.Item(19,5).Value2 = "  -3.80%  "

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value2 = "420.78"
$ws.Cells.Item(20,5).Value2 = "  -0.72%  "

# Row 21
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value2 = "13.22"
$ws.Cells.Item(21,5).Value2 = "  +1.60%  "

# Row 22
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value2 = "0.686"
$ws.Cells.Item(22,5).Value2 = "  +4.41%  "

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value2 = "7.04"
$ws.Cells.Item(23,5).Value2 = "  -0.63%  "

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value2 = "13.04"
$ws.Cells.Item(24,5).Value2 = "  +1.71%  "

# Row 25
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value2 = "79.66"
$ws.Cells.Item(25,5).Value2 = "  +1.43%  "

# Row 26
$ws.Cells.Item(26,5).Value2 = "  +0.21%  "

# Row 27
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value2 = "0.995"
$ws.Cells.Item(27,5).Value2 = "  -0.45%  "

# Row 28
$ws.Cells.Item(28,5).Value2 = "  +0.23%  "

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value2 = "7.65"
$ws.Cells.Item(29,5).Value2 = "  +7.31%  "

# Row 30
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value2 = "2.00"
$ws.Cells.Item(30,5).Value2 = "  +6.86%  "

# Row 31
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value2 = "6.21"
$ws.Cells.Item(31,5).Value2 = "  +2.78%  "

# Row 32
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value2 = "0.105"
$ws.Cells.Item(32,5).Value2 = "  +13.79%  "

# Row 33
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value2 = "25.32"
$ws.Cells.Item(33,5).Value2 = "  +0.99%  "

# Row 34
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value2 = "5.66"
$ws.Cells.Item(34,5).Value2 = "  +2.39%  "

# Row 35
$ws.Cells.Item(35,2).Value2 = "Stacks"
$ws.Cells.Item(35,3).Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value2 = "2.14"
$ws.Cells.Item(35,5).Value2 = "  -4.24%  "

# Row 36
$ws.Cells.Item(36,2).Value2 = "Mantle"
$ws.Cells.Item(36,3).Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value2 = "0.942"
$ws.Cells.Item(36,5).Value2 = "  -0.37%  "

# Row 37
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value2 = "48.58"
$ws.Cells.Item(37,5).Value2 = "  -1.51%  "

# Row 38
$ws.Cells.Item(38,4).Value2 = "0.0₃0672"
$ws.Cells.Item(38,5).Value2 = "  +3.83%  "

# Row 39
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value2 = "8.38"
$ws.Cells.Item(39,5).Value2 = "  +6.95%  "

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value2 = "2.62"
$ws.Cells.Item(40,5).Value2 = "  +10.50%  "

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value2 = "0.109"
$ws.Cells.Item(41,5).Value2 = "  +0.94%  "

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value2 = "0.0349"
$ws.Cells.Item(42,5).Value2 = "  -1.76%  "

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value2 = "378.95"
$ws.Cells.Item(43,5).Value2 = "  +0.86%  "

# Row 44
$ws.Cells.Item(44,4).Value2 = "2.661.61"
$ws.Cells.Item(44,5).Value2 = "  +1.18%  "

# Row 46
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value2 = "0.241"
$ws.Cells.Item(46,5).Value2 = "  +3.29%  "

# Row 47
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value2 = "120.57"
$ws.Cells.Item(47,5).Value2 = "  +1.52%  "

# Row 48
$ws.Cells.Item(48,5).Value2 = "  +3.75%  "

# Row 49
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value2 = "2.01"
$ws.Cells.Item(49,5).Value2 = "  +2.84%  "

# Row 50
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value2 = "23.60"
$ws.Cells.Item(50,5).Value2 = "  +1.71%  "

# Row 51
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value2 = "2.01"
$ws.Cells.Item(51,5).Value2 = "  +2.42%  "
